$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 6 data: SerialNo=5, Topic=Insertion Sort, Code/Algo=c, T.C=O(n2), S.C=O(1), ClassName=InsertionSort
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Insertion Sort"
$ws.Range("D6").Value = "c"
$ws.Range("E6").Value = "O(n2)"
$ws.Range("F6").Value = "O(1)"
$ws.Range("H6").Value = "InsertionSort"

# Update selection to match new active cell
$ws.Range("H6").Select()

$wb.Save()
